$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) First paragraph: "This is a Microsoft word document." gets two
#    trailing spaces appended, then three new red (FF0000) runs are
#    appended before the paragraph mark:
#      "(This is a change – Ve" + "rsion for main branch" + ")"
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1)

# Append the two trailing spaces to the existing (black) run by
# inserting right before the paragraph mark.
$insPoint = $p1.Range.End - 1
$spaceRange = $d.Range($insPoint, $insPoint)
$spaceRange.InsertAfter("  ")

# Run 2: "(This is a change – Ve"
$insPoint = $p1.Range.End - 1
$run2 = $d.Range($insPoint, $insPoint)
$run2.InsertAfter([string][char]40 + "This is a change " + [string][char]8211 + " Ve")
$run2Fmt = $d.Range($insPoint, $p1.Range.End - 1)
$run2Fmt.Font.Color = 255

# Run 3: "rsion for main branch"
$insPoint = $p1.Range.End - 1
$run3 = $d.Range($insPoint, $insPoint)
$run3.InsertAfter("rsion for main branch")
$run3Fmt = $d.Range($insPoint, $p1.Range.End - 1)
$run3Fmt.Font.Color = 255

# Run 4: ")"
$insPoint = $p1.Range.End - 1
$run4 = $d.Range($insPoint, $insPoint)
$run4.InsertAfter([string][char]41)
$run4Fmt = $d.Range($insPoint, $p1.Range.End - 1)
$run4Fmt.Font.Color = 255

# ---------------------------------------------------------------------
# 2) Drop the trailing "... ank God almighty, we are free at last."
#    paragraph entirely (it was the final paragraph of the document).
# ---------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$lastPara.Range.Delete()
